$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '96.431.33'
$ws.Range('E2').Value = '  -0.76%  '
$ws.Range('D3').Value = '3.634.34'
$ws.Range('E3').Value = '  -2.29%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '2.71'
$ws.Range('E4').Value = '  +42.34%  '
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '226.56'
$ws.Range('E6').Value = '  -4.62%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '642.65'
$ws.Range('E7').Value = '  -2.60%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.426'
$ws.Range('E8').Value = '  -2.02%  '
$ws.Range('E9').Value = '  +14.53%  '
$ws.Range('E10').Value = '  +0.03%  '
$ws.Range('D11').Value = '3.633.44'
$ws.Range('E11').Value = '  -2.28%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '48.51'
$ws.Range('E12').Value = '  +8.11%  '
$ws.Range('E13').Value = '  +2.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000292'
$ws.Range('E14').Value = '  -8.47%  '
$ws.Range('E15').Value = '  -5.67%  '
$ws.Range('D16').Value = '4.313.66'
$ws.Range('E16').Value = '  -2.34%  '
$ws.Range('D17').Value = '96.239.36'
$ws.Range('E17').Value = '  -0.81%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '22.00'
$ws.Range('E18').Value = '  +16.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '8.90'
$ws.Range('E19').Value = '  -1.26%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.98'
$ws.Range('E20').Value = '  +6.92%  '
$ws.Range('D21').Value = '3.636.68'
$ws.Range('E21').Value = '  -2.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.557'
$ws.Range('E22').Value = '  +10.31%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.280'
$ws.Range('E23').Value = '  +46.62%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '518.06'
$ws.Range('E24').Value = '  -1.57%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.26'
$ws.Range('E25').Value = '  -6.23%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '121.87'
$ws.Range('E26').Value = '  +13.99%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000202'
$ws.Range('E27').Value = '  -10.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.80'
$ws.Range('E28').Value = '  -1.34%  '
$ws.Range('D29').Value = '3.816.80'
$ws.Range('E29').Value = '  -2.99%  '
$ws.Range('E30').Value = '  -5.61%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '12.96'
$ws.Range('E31').Value = '  +1.43%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.04'
$ws.Range('E32').Value = '  -0.11%  '
$ws.Range('E33').Value = '  +0.10%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.622'
$ws.Range('E34').Value = '  +5.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '32.76'
$ws.Range('E35').Value = '  +0.59%  '
$ws.Range('B36').Value = 'Cronos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.180'
$ws.Range('E36').Value = '  -6.69%  '
$ws.Range('B37').Value = 'Binance-PegBSC-USD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  -0.11%  '
$ws.Range('E38').Value = '  -4.35%  '
$ws.Range('E39').Value = '  -0.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.40'
$ws.Range('E40').Value = '  -4.18%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '583.10'
$ws.Range('E41').Value = '  -9.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '7.09'
$ws.Range('E42').Value = '  +5.21%  '
$ws.Range('E43').Value = '  +4.55%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0510'
$ws.Range('E44').Value = '  +11.74%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '40.85'
$ws.Range('E45').Value = '  +0.79%  '
$ws.Range('E46').Value = '  -5.36%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.957'
$ws.Range('E47').Value = '  -1.49%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.95'
$ws.Range('E48').Value = '  -3.69%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '230.84'
$ws.Range('E49').Value = '  +12.15%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.90'
$ws.Range('E50').Value = '  +3.00%  '
$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '23.59'
$ws.Range('E51').Value = '  -0.20%  '
